$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.990.00'
$ws.Range("E2").Value = '  +5.06%  '
$ws.Range("D3").Value = '2.335.59'
$ws.Range("E3").Value = '  +2.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.60'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.50'
$ws.Range("E6").Value = '  +4.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").Value = '  +2.95%  '
$ws.Range("D9").Value = '2.369.88'
$ws.Range("E9").Value = '  +4.02%  '
$ws.Range("E10").Value = '  +9.39%  '
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.21'
$ws.Range("E12").Value = '  +6.25%  '
$ws.Range("E13").Value = '  +3.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.10'
$ws.Range("E14").Value = '  +4.38%  '
$ws.Range("D15").Value = '2.756.72'
$ws.Range("E15").Value = '  +3.19%  '
$ws.Range("D16").Value = '57.350.04'
$ws.Range("E16").Value = '  +5.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  +5.30%  '
$ws.Range("D18").Value = '2.344.17'
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.64'
$ws.Range("E19").Value = '  +3.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.31'
$ws.Range("E20").Value = '  +3.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.82'
$ws.Range("E21").Value = '  +6.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.68'
$ws.Range("E22").Value = '  +5.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.76'
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.991'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  +7.43%  '
$ws.Range("E27").Value = '  +6.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.34'
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").Value = '0.0₃0752'
$ws.Range("E29").Value = '  +6.87%  '
$ws.Range("E30").Value = '  +12.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.36'
$ws.Range("E31").Value = '  +6.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.70'
$ws.Range("E32").Value = '  +5.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.48'
$ws.Range("E33").Value = '  +3.70%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.961'
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("E36").Value = '  +5.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.991'
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("E38").Value = '  +9.08%  '
$ws.Range("E39").Value = '  +9.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.66'
$ws.Range("E40").Value = '  +4.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.383'
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.20'
$ws.Range("E42").Value = '  +11.92%  '
$ws.Range("E43").Value = '  +7.46%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '279.45'
$ws.Range("E44").Value = '  +14.32%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.16'
$ws.Range("E45").Value = '  +5.00%  '
$ws.Range("E46").Value = '  +4.11%  '
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.569'
$ws.Range("E48").Value = '  +4.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0218'
$ws.Range("E49").Value = '  +6.23%  '
$ws.Range("E50").Value = '  +2.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.09'
$ws.Range("E51").Value = '  +5.39%  '
